$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 1, shifting all existing data down.
$ws.Rows.Item(1).Insert()

# Add header labels for the data columns (imported as data point labels).
$ws.Range("A1").Value = "t"
$ws.Range("B1").Value = "One"
$ws.Range("C1").Value = "Two"
$ws.Range("D1").Value = "Three"

# Match the saved selection state (active cell on the new header row).
$ws.Range("D1").Select()
